$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.786.44'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.85%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.428.87'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.03%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '552.94'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.94%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '160.21'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.70%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.496'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.68%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.426.71'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.07%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.146'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -6.80%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.163'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.86%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.331'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -5.67%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.70'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.69%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.875.95'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '67.767.09'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.81%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000164'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -5.64%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '22.76'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -5.73%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.464.36'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.67'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.93%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '336.39'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.35%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.97'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -5.16%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.68'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.97%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.81'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.74%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '65.78'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -5.05%  '
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.58'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -7.00%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.01'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.95'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -7.69%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0₃0808'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -6.56%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.99'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -8.46%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.998'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '412.10'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -6.18%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.61'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.93%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.10'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -6.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '157.19'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.28%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.96'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.106'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.95%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '17.59'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.57%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.297'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.69%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.24'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -6.89%  '
$ws.Range('E43').Value = '  -7.32%  '
$ws.Range('E44').Value = '  -0.84%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '131.86'
$ws.Range('D45').Style = "Normal"
$ws.Range('E46').Value = '  -6.76%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.27'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.26%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0709'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.87%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.467'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -8.16%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.550'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.24%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0898'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.08%  '
